$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as literal text
# (so number-looking strings like "1.001" keep their exact text form),
# without leaving a residual numeric cell style behind.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "28.415.57"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.822.42"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "314.33"
$ws.Range("E5").Value = "  -0.73%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.5125"
$ws.Range("E7").Value = "  -3.37%  "
Set-TextValue $ws.Range("D8") "0.3930"
$ws.Range("E8").Value = "  -3.36%  "
Set-TextValue $ws.Range("D9") "0.07674"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E12").Value = "  +0.62%  "
Set-TextValue $ws.Range("D13") "6.269"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  +0.07%  "
Set-TextValue $ws.Range("D15") "7.499"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "1.821.62"
$ws.Range("E16").Value = "  -0.82%  "
Set-TextValue $ws.Range("D17") "93.00"
$ws.Range("E17").Value = "  +3.84%  "
Set-TextValue $ws.Range("D18") "0.00001105"
$ws.Range("E18").Value = "  +2.91%  "
Set-TextValue $ws.Range("D19") "0.06646"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  -0.01%  "
Set-TextValue $ws.Range("D22") "6.099"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "28.432.09"
$ws.Range("E23").Value = "  -0.25%  "
Set-TextValue $ws.Range("D24") "11.17"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  +5.67%  "
Set-TextValue $ws.Range("D26") "20.86"
$ws.Range("E26").Value = "  +1.38%  "
Set-TextValue $ws.Range("D27") "156.26"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "2.033.47"
$ws.Range("E28").Value = "  -0.58%  "
Set-TextValue $ws.Range("D29") "2.388"
$ws.Range("E29").Value = "  -3.31%  "
Set-TextValue $ws.Range("D30") "123.97"
$ws.Range("E30").Value = "  +0.24%  "
Set-TextValue $ws.Range("D31") "0.1100"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").Value = "  -1.28%  "
Set-TextValue $ws.Range("D33") "5.643"
$ws.Range("E33").Value = "  -0.76%  "
Set-TextValue $ws.Range("D34") "3.655"
$ws.Range("E34").Value = "  +0.01%  "
Set-TextValue $ws.Range("D35") "0.07074"
$ws.Range("E35").Value = "  -1.00%  "
Set-TextValue $ws.Range("D36") "0.2210"
$ws.Range("E36").Value = "  -2.65%  "
Set-TextValue $ws.Range("D37") "0.02328"
$ws.Range("E37").Value = "  -0.71%  "
Set-TextValue $ws.Range("D38") "5.165"
$ws.Range("E38").Value = "  -1.87%  "
Set-TextValue $ws.Range("D39") "8.744"
$ws.Range("E39").Value = "  -0.71%  "
Set-TextValue $ws.Range("D40") "0.6256"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -1.39%  "
Set-TextValue $ws.Range("D42") "1.171"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("E43").Value = "  -0.03%  "
Set-TextValue $ws.Range("D44") "1.391"
$ws.Range("E44").Value = "  -1.19%  "
Set-TextValue $ws.Range("D45") "13.34"
$ws.Range("E45").Value = "  -0.98%  "
Set-TextValue $ws.Range("D46") "3.728"
$ws.Range("E46").Value = "  +0.59%  "
Set-TextValue $ws.Range("D47") "0.5872"
$ws.Range("E47").Value = "  +0.21%  "
Set-TextValue $ws.Range("D48") "124.11"
$ws.Range("E48").Value = "  -1.58%  "
Set-TextValue $ws.Range("D49") "1.978"
$ws.Range("E49").Value = "  -0.71%  "
Set-TextValue $ws.Range("D50") "1.193"
Set-TextValue $ws.Range("D51") "0.06888"
$ws.Range("E51").Value = "  -0.19%  "
